$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The first line item (row 8: "Rewiring of light point...") is being dropped from
# the bill. Delete its row -- this shifts every following line item (and the
# summary rows below) up by one, which already gives every remaining row the
# correct label/description/rate cells.
$ws.Rows(8).Delete()

# After the shift, two of the "P. point" rows no longer carry that unit label.
$ws.Range("A9").Value = ""
$ws.Range("A10").Value = ""

# Updated measured quantities ("Qty executed upto date") for the remaining
# line items.
$ws.Range("C8").Value = 82
$ws.Range("C9").Value = 79
$ws.Range("C10").Value = 20
$ws.Range("C11").Value = 34
$ws.Range("C12").Value = 52
$ws.Range("C13").Value = 39
$ws.Range("C14").Value = 13
$ws.Range("C15").Value = 10
$ws.Range("C16").Value = 17

# Recomputed "Upto date Amount" (Qty * Rate) for the rows whose quantity
# changed above. These are kept as text (matching the sheet's existing
# "0.00"-style values) -- briefly force Text format so the numeric-looking
# string isn't auto-converted to a number, then restore the Normal style so
# the cell's formatting is left exactly as it was.
function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "G8" "20992.00"
Set-TextValue "G9" "37288.00"
Set-TextValue "G10" "13240.00"
Set-TextValue "G12" "7072.00"
Set-TextValue "G13" "897.00"

# Grand total / net payable amount rows pick up the new overall total.
Set-TextValue "G18" "79489.00"
Set-TextValue "H18" "79489.00"
Set-TextValue "G20" "79489.00"
Set-TextValue "H20" "79489.00"
